# ---------------------------------------------------------------------------
# Edit script: mirrors the OOXML diff that
#   1) re-points the three Table_0-styled tables (slides 14, 15, 16) at the
#      built-in table style {387B8FD6-0D87-45DE-9FAB-C6CF8FDF7E52}
#      (previously {C459DA5E-F46A-4DE1-BB8C-CC6BDD0C6726}), and
#   2) swaps the deck's applied theme from the custom "Integral" (Red Violet)
#      design back to the default "Office Theme" colors.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Re-apply the (built-in) table style to the three tables that currently
#    carry the custom "Table_0" style id.
# ---------------------------------------------------------------------------
$newTableStyleId = "{387B8FD6-0D87-45DE-9FAB-C6CF8FDF7E52}"
$tableSlideIndexes = @(14, 15, 16)

foreach ($slideIdx in $tableSlideIndexes) {
    $slide = $p.Slides.Item($slideIdx)
    $shape = $slide.Shapes.Item(1)
    if ($shape.HasTable) {
        $table = $shape.Table
        $table.ApplyStyle($newTableStyleId)
    }
}

# ---------------------------------------------------------------------------
# 2) Apply the default "Office Theme" color palette across the deck (the
#    design switch that, on save, leaves the Office Theme colors/fonts in
#    the theme part the slide master/slides actually use).
#    Theme color order (matches MsoThemeColorSchemeIndex / the XML
#    <a:clrScheme> child order): dk1, lt1, dk2, lt2, accent1-6, hlink,
#    folHlink.
# ---------------------------------------------------------------------------
function Hex-To-RGB($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$officeThemeColors = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

$slide1 = $p.Slides.Item(1)
$themeColors = $slide1.ThemeColorScheme

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $entry = $themeColors.Item($i)
    $entry.RGB = Hex-To-RGB $officeThemeColors[$i - 1]
}
